$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) New bullet list item after "treeman: tree manipulation"
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "adephylo: compute distance from the root to the tips (distRoot)"

# ---------------------------------------------------------------------------
# 2) Normal style tweaks: font color, bidi / justification
# ---------------------------------------------------------------------------
$normal = $d.Styles.Item("Normal")
$normal.Font.Color = 655360          # RGB(0x0A,0x00,0x00) -> w:color val="00000A"
$normal.ParagraphFormat.ReadingOrder = 2   # -> <w:bidi w:val="0"/>
$normal.ParagraphFormat.Alignment = 0      # -> <w:jc w:val="left"/>

# ---------------------------------------------------------------------------
# 3) New character styles ListLabel1 .. ListLabel9
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 9; $i++) {
    $style = $d.Styles.Add("ListLabel $i", 2)
    $style.QuickStyle = $true
    $style.Font.NameBi = "OpenSymbol"
}
